# Clear out the stale "Tarea c" task row (row 4) and drop the now-unused
# trailing blank row (row 10), matching the cleanup done while wiring up
# dbxls / functions / menuapp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the task data that used to live in row 4 (Id, Titulo, Descripcion, Estado).
$ws.Range("A4:D4").ClearContents()

# Keep row 4 itself present (as a blank row spanning A:F) instead of letting it
# collapse away entirely, and make sure none of its cells pick up formatting.
$ws.Range("A4:F4").Style = "Normal"

# The worksheet no longer needs the empty trailing row 10.
$ws.Rows("10").Delete()
